$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(106, 8).Value = 3479.375
$ws.Cells.Item(106, 9).Value = 3262.1428
$ws.Cells.Item(106, 11).Value = 3262.1428
$ws.Cells.Item(106, 13).Value = -2631.1428

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4489.839
$ws.Cells.Item(32, 9).Value = 3733.673
$ws.Cells.Item(32, 11).Value = 3733.673
$ws.Cells.Item(32, 13).Value = -3446.673

$ws.Cells.Item(43, 8).Value = 31757
$ws.Cells.Item(43, 10).Value = 31757
$ws.Cells.Item(43, 12).Value = 31757
$ws.Cells.Item(43, 14).Value = -32383

$ws.Cells.Item(45, 8).Value = 6853487.5
$ws.Cells.Item(45, 9).Value = 13079305
$ws.Cells.Item(45, 11).Value = 13079305
$ws.Cells.Item(45, 13).Value = -13078928

$ws.Cells.Item(110, 8).Value = 1639008.8
$ws.Cells.Item(110, 9).Value = 1740996.8
$ws.Cells.Item(110, 11).Value = 1740996.8
$ws.Cells.Item(110, 13).Value = -1738951.8

$ws.Cells.Item(122, 8).Value = 497969.1
$ws.Cells.Item(122, 9).Value = 1582.0555
$ws.Cells.Item(122, 11).Value = 4746.166499999999
$ws.Cells.Item(122, 13).Value = -2296.166499999999

$ws.Cells.Item(132, 8).Value = 2464.3
$ws.Cells.Item(132, 9).Value = 1812.8889
$ws.Cells.Item(132, 10).Value = 3441.4167
$ws.Cells.Item(132, 11).Value = 5438.6667
$ws.Cells.Item(132, 12).Value = 10324.2501
$ws.Cells.Item(132, 13).Value = -2908.6667
$ws.Cells.Item(132, 14).Value = -15384.2501

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(95, 8).Value = 9156
$ws.Cells.Item(95, 10).Value = 9156
$ws.Cells.Item(95, 12).Value = 9156
$ws.Cells.Item(95, 14).Value = -14648

$ws.Cells.Item(107, 8).Value = 6497048.5
$ws.Cells.Item(107, 9).Value = 10204933
$ws.Cells.Item(107, 11).Value = 10204933
$ws.Cells.Item(107, 13).Value = -10203013

$ws.Cells.Item(134, 8).Value = 2599.8823
$ws.Cells.Item(134, 9).Value = 1062.8718
$ws.Cells.Item(134, 10).Value = 7595.1665
$ws.Cells.Item(134, 11).Value = 3188.6154
$ws.Cells.Item(134, 12).Value = 22785.4995
$ws.Cells.Item(134, 13).Value = -653.6153999999997
$ws.Cells.Item(134, 14).Value = -27855.4995

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 30328.438
$ws.Cells.Item(31, 10).Value = 85410.17999999999
$ws.Cells.Item(31, 12).Value = 85410.17999999999
$ws.Cells.Item(31, 14).Value = -86000.17999999999

$ws.Cells.Item(34, 8).Value = 30328.438
$ws.Cells.Item(34, 10).Value = 85410.17999999999
$ws.Cells.Item(34, 12).Value = 85410.17999999999
$ws.Cells.Item(34, 14).Value = -85814.17999999999

$ws.Cells.Item(58, 8).Value = 2002.2413
$ws.Cells.Item(58, 9).Value = 1648.6666
$ws.Cells.Item(58, 10).Value = 2580.818
$ws.Cells.Item(58, 11).Value = 1648.6666
$ws.Cells.Item(58, 12).Value = 2580.818
$ws.Cells.Item(58, 13).Value = -1445.6666
$ws.Cells.Item(58, 14).Value = -2986.818

$ws.Cells.Item(94, 8).Value = 1333.1111
$ws.Cells.Item(94, 9).Value = 1199
$ws.Cells.Item(94, 11).Value = 1199
$ws.Cells.Item(94, 13).Value = -748

$ws.Cells.Item(122, 8).Value = 1880.1765
$ws.Cells.Item(122, 9).Value = 1464.8667
$ws.Cells.Item(122, 11).Value = 4394.6001
$ws.Cells.Item(122, 13).Value = -1944.6001

$ws.Cells.Item(132, 8).Value = 80291.164
$ws.Cells.Item(132, 9).Value = 49114.523
$ws.Cells.Item(132, 11).Value = 147343.569
$ws.Cells.Item(132, 13).Value = -144813.569

$ws.Cells.Item(134, 8).Value = 36909.742
$ws.Cells.Item(134, 10).Value = 4579.2144
$ws.Cells.Item(134, 12).Value = 13737.6432
$ws.Cells.Item(134, 14).Value = -18807.6432

$ws.Cells.Item(136, 8).Value = 2002.2413
$ws.Cells.Item(136, 9).Value = 1648.6666
$ws.Cells.Item(136, 10).Value = 2580.818
$ws.Cells.Item(136, 11).Value = 4945.9998
$ws.Cells.Item(136, 12).Value = 7742.454000000001
$ws.Cells.Item(136, 13).Value = -2395.9998
$ws.Cells.Item(136, 14).Value = -12842.454

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(56, 8).Value = 9620519
$ws.Cells.Item(56, 9).Value = 9620519
$ws.Cells.Item(56, 11).Value = 9620519
$ws.Cells.Item(56, 13).Value = -9619989

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(53, 8).Value = 50000
$ws.Cells.Item(53, 10).Value = 50000
$ws.Cells.Item(53, 12).Value = 50000
$ws.Cells.Item(53, 14).Value = -51262

$ws.Cells.Item(70, 8).Value = 16673998
$ws.Cells.Item(70, 9).Value = 20007598
$ws.Cells.Item(70, 11).Value = 20007598
$ws.Cells.Item(70, 13).Value = -20007328

$ws.Cells.Item(73, 8).Value = 16673998
$ws.Cells.Item(73, 9).Value = 20007598
$ws.Cells.Item(73, 11).Value = 20007598
$ws.Cells.Item(73, 13).Value = -20006662

$ws.Cells.Item(102, 8).Value = 6476048
$ws.Cells.Item(102, 9).Value = 8550192
$ws.Cells.Item(102, 11).Value = 8550192
$ws.Cells.Item(102, 13).Value = -8548570

$ws.Cells.Item(113, 8).Value = 13890570
$ws.Cells.Item(113, 9).Value = 16668289
$ws.Cells.Item(113, 10).Value = 1975
$ws.Cells.Item(113, 11).Value = 16668289
$ws.Cells.Item(113, 12).Value = 1975
$ws.Cells.Item(113, 13).Value = -16666119
$ws.Cells.Item(113, 14).Value = -6315

$ws.Cells.Item(122, 8).Value = 389206.53
$ws.Cells.Item(122, 10).Value = 3244.75
$ws.Cells.Item(122, 12).Value = 9734.25
$ws.Cells.Item(122, 14).Value = -14634.25

$ws.Cells.Item(127, 8).Value = 50126
$ws.Cells.Item(127, 10).Value = 50126
$ws.Cells.Item(127, 12).Value = 50126
$ws.Cells.Item(127, 14).Value = -60046

$ws.Cells.Item(132, 8).Value = 3481.4211
$ws.Cells.Item(132, 9).Value = 2842.7856
$ws.Cells.Item(132, 10).Value = 5269.6
$ws.Cells.Item(132, 11).Value = 8528.356800000001
$ws.Cells.Item(132, 12).Value = 15808.8
$ws.Cells.Item(132, 13).Value = -5998.356800000001
$ws.Cells.Item(132, 14).Value = -20868.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 4586.9106
$ws.Cells.Item(132, 9).Value = 3969.3428
$ws.Cells.Item(132, 10).Value = 5616.1904
$ws.Cells.Item(132, 11).Value = 11908.0284
$ws.Cells.Item(132, 12).Value = 16848.5712
$ws.Cells.Item(132, 13).Value = -9378.028399999999
$ws.Cells.Item(132, 14).Value = -21908.5712

$ws.Cells.Item(136, 8).Value = 114763.39
$ws.Cells.Item(136, 9).Value = 146671.58
$ws.Cells.Item(136, 10).Value = 3084.75
$ws.Cells.Item(136, 11).Value = 440014.74
$ws.Cells.Item(136, 12).Value = 9254.25
$ws.Cells.Item(136, 13).Value = -437464.74
$ws.Cells.Item(136, 14).Value = -14354.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(23, 8).Value = 463.33334
$ws.Cells.Item(23, 9).Value = 463.33334
$ws.Cells.Item(23, 11).Value = 463.33334
$ws.Cells.Item(23, 13).Value = -234.33334

$ws.Cells.Item(30, 8).Value = 18890
$ws.Cells.Item(30, 10).Value = 18890
$ws.Cells.Item(30, 12).Value = 18890
$ws.Cells.Item(30, 14).Value = -19104

$ws.Cells.Item(100, 8).Value = 1284.9375
$ws.Cells.Item(100, 10).Value = 4496.3335
$ws.Cells.Item(100, 12).Value = 8992.666999999999
$ws.Cells.Item(100, 14).Value = -10074.667

$ws.Cells.Item(107, 8).Value = 41668980

$ws.Cells.Item(129, 8).Value = 39499
$ws.Cells.Item(129, 10).Value = 39499
$ws.Cells.Item(129, 12).Value = 39499
$ws.Cells.Item(129, 14).Value = -49499

$ws.Cells.Item(132, 8).Value = 45952210
$ws.Cells.Item(132, 9).Value = 58825596
$ws.Cells.Item(132, 11).Value = 176476788
$ws.Cells.Item(132, 13).Value = -176474258
